$d = $word.ActiveDocument

$newText = "Sixthly, we found linux platform could be used as the operating system if we use virtualization tools like virtualbox. In out case, we installed Ubuntu guest os on the windows virtualbox.  In order to use the serial port on the guest machine, we need to map the serial port in the host machine to guest machine. Then we found when the serial port under windows[host] won’t work as long as the linux[guest] is working . Another issue is the /dev/ttyS0 device will be generated automatically after booting , however, it requires root access to manipulate other wise the port will reject any requests from our program."

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last
$newLast.Range.Text = $newText
